$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restated financial figures (unit-scale correction) for rows 2-9,
# columns D:AJ. Built as a row => { col = value } map and applied in one pass.
$rowData = @{
    2 = @{ "D" = 13839; "E" = 812; "F" = 812; "G" = 777; "H" = 586; "I" = 607; "J" = -20; "K" = 14649; "L" = 7754; "M" = 6895; "N" = 6904; "O" = -9; "P" = 168; "Q" = 1180; "R" = -1709; "S" = 684; "T" = 984; "U" = 197; "V" = 5795; "W" = 5.87; "X" = 4.24; "Y" = 9.18; "Z" = 4.28; "AA" = 112.46; "AB" = 4103.47; "AC" = 18034; "AD" = 17.24; "AE" = 205268; "AF" = 1.52; "AG" = 2000; "AH" = 0.64; "AI" = 11.09; "AJ" = 3363200 }
    3 = @{ "D" = 13597; "E" = 573; "F" = 573; "G" = 211; "H" = 70; "I" = 116; "J" = -46; "K" = 16180; "L" = 9091; "M" = 7089; "N" = 7129; "O" = -40; "P" = 168; "Q" = 560; "R" = -575; "S" = 732; "T" = 875; "U" = -316; "V" = 6913; "W" = 4.22; "X" = 0.52; "Y" = 1.65; "Z" = 0.46; "AA" = 128.24; "AB" = 4132.49; "AC" = 3451; "AD" = 76.64; "AE" = 211974; "AF" = 1.25; "AG" = 2000; "AH" = 0.76; "AI" = 57.95; "AJ" = 3363200 }
    4 = @{ "D" = 15764; "E" = 1514; "F" = 1514; "G" = 1469; "H" = 774; "I" = 750; "J" = 24; "K" = 19003; "L" = 11092; "M" = 7911; "N" = 7925; "O" = -14; "P" = 168; "Q" = 2101; "R" = -1183; "S" = 945; "T" = 729; "U" = 1373; "V" = 8208; "W" = 9.609999999999999; "X" = 4.91; "Y" = 9.960000000000001; "Z" = 4.4; "AA" = 140.21; "AB" = 4538.49; "AC" = 22300; "AD" = 16.3; "AE" = 235636; "AF" = 1.54; "AG" = 3000; "AH" = 0.83; "AI" = 13.45; "AJ" = 3363200 }
    5 = @{ "D" = 23811; "E" = 2246; "F" = 2246; "G" = 2202; "H" = 1678; "I" = 1668; "J" = 10; "K" = 25603; "L" = 16566; "M" = 9037; "N" = 9061; "O" = -23; "P" = 168; "Q" = 2811; "R" = -4175; "S" = -169; "T" = 1546; "U" = 1265; "V" = 11389; "W" = 9.43; "X" = 7.05; "Y" = 19.65; "Z" = 7.52; "AA" = 183.31; "AB" = 5467.56; "AC" = 49609; "AD" = 6.77; "AE" = 269402; "AF" = 1.25; "AG" = 4000; "AH" = 1.19; "AI" = 8.06; "AJ" = 3363200 }
    6 = @{ "D" = 24446; "E" = 1777; "F" = 1777; "G" = 903; "H" = 710; "I" = 691; "K" = 25922; "L" = 16186; "M" = 9736; "N" = 9750; "P" = 168; "Q" = 528; "R" = -663; "S" = -898; "T" = 707; "U" = -179; "V" = 11080; "W" = 7.27; "X" = 2.9; "Y" = 7.35; "Z" = 2.76; "AA" = 166.25; "AB" = 5769.37; "AC" = 20549; "AD" = 9.880000000000001; "AE" = 289909; "AF" = 0.7; "AG" = 3000; "AH" = 1.48; "AI" = 14.6; "AJ" = 3363200 }
    7 = @{ "D" = 27260; "E" = 2230; "G" = 1220; "H" = 880; "I" = 730; "K" = 32690; "L" = 20780; "M" = 11910; "N" = 11440; "P" = 180; "U" = 1170; "W" = 8.18; "X" = 3.23; "Y" = 6.89; "Z" = 3; "AA" = 174.48; "AC" = 20573; "AD" = 9.82; "AE" = 311069; "AF" = 0.65; "AG" = 3500; "AH" = 1.73; "AI" = 17.63 }
    8 = @{ "D" = 29300; "E" = 2350; "G" = 1740; "H" = 1270; "I" = 1230; "K" = 33980; "L" = 20930; "M" = 13050; "N" = 12540; "P" = 180; "U" = 1480; "W" = 8.02; "X" = 4.33; "Y" = 10.26; "Z" = 3.81; "AA" = 160.38; "AC" = 33445; "AD" = 6.04; "AE" = 340979; "AF" = 0.59; "AG" = 3500; "AH" = 1.73; "AI" = 10.46 }
    9 = @{ "D" = 31200; "E" = 2500; "G" = 1920; "H" = 1380; "I" = 1340; "K" = 35620; "L" = 21310; "M" = 14300; "N" = 13750; "P" = 180; "U" = 1590; "W" = 8.01; "X" = 4.42; "Y" = 10.19; "Z" = 3.97; "AA" = 149.02; "AC" = 36436; "AD" = 5.54; "AE" = 373881; "AF" = 0.54; "AG" = 3500; "AH" = 1.73; "AI" = 9.609999999999999 }
}

foreach ($rowNum in $rowData.Keys) {
    $cols = $rowData[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}

# Rows 7-9 dropped their Q/R/S/T figures entirely (columns no longer reported)
foreach ($rowNum in 7,8,9) {
    foreach ($col in "Q","R","S","T") {
        $ws.Range("$col$rowNum").ClearContents()
    }
}